# Applies the "updated output" changes described in the commit
# "Added calcium and b12 to Zambia and Uganda and updated output" for
# output/zinc_SPADE_uganda_h/3_excel_tables/spade_uganda_h_zinc.xlsx
#
# Changes:
#  - Info sheet: refresh the run Start_time / End_time timestamps.
#  - sessionInfo sheet: bump package version numbers that changed between
#    the two runs (here, magrittr, rprojroot).
#  - sessionInfo sheet: the "Loaded_only" package table lost one entry
#    (backports) compared to the previous run, so the last populated row
#    (which held "boot" / "1.3-25") becomes empty.

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Info")
$wsSession = $wb.Worksheets.Item("sessionInfo")

# --- Info sheet: Start_time / End_time (column B, rows 26/27) ---
$wsInfo.Range("B26").Value = "Thu Nov 19 15:23:55 2020"
$wsInfo.Range("B27").Value = "Thu Nov 19 15:23:59 2020"

# --- sessionInfo sheet: updated "Ohter_packages" version numbers ---
# here: 0.1 -> 1.0.0
$wsSession.Range("G2").Value = "1.0.0"

# --- sessionInfo sheet: updated "Loaded_only" version numbers ---
# magrittr: 1.5 -> 2.0.1
$wsSession.Range("J3").Value = "2.0.1"
# rprojroot: 1.3-2 -> 2.0.2
$wsSession.Range("J10").Value = "2.0.2"

# --- sessionInfo sheet: "backports" package no longer loaded, so the
# "Loaded_only" list (I15:J15 = backports/1.1.10) is removed and the rows
# below it shift up by one; the previous last row (16), which used to hold
# "boot" / "1.3-25", becomes blank.
$wsSession.Range("I15").Value = $wsSession.Range("I16").Value2
$wsSession.Range("J15").Value = $wsSession.Range("J16").Value2
$wsSession.Range("I16:J16").ClearContents()
